# Update gh-pages to output generated at 456a3b4
# This script applies the numeric "F" column updates (view/count style counters)
# and one event-title edit (removal of a "mystery guest" suffix) across the
# four worksheets of the 北京-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 309
$ws1.Range("F4").Value  = 436
$ws1.Range("F5").Value  = 8636
$ws1.Range("F7").Value  = 10878
$ws1.Range("F22").Value = 1844
$ws1.Range("F24").Value = 589
$ws1.Range("F30").Value = 1228
$ws1.Range("F35").Value = 1422
$ws1.Range("F37").Value = 349
$ws1.Range("F39").Value = 28
$ws1.Range("F41").Value = 524
$ws1.Range("F42").Value = 358
$ws1.Range("F43").Value = 105
$ws1.Range("F44").Value = 806

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 48
$ws2.Range("C14").Value = "北京·法国姐姐”乔伊丝·乔纳森《小意思》巡回演唱会"

# ---- Sheet "本地生活" (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 214
$ws3.Range("F3").Value = 2820
$ws3.Range("F5").Value = 210

# ---- Sheet "全部类型" (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 214
$ws4.Range("F6").Value  = 210
$ws4.Range("F8").Value  = 436
$ws4.Range("F9").Value  = 8636
$ws4.Range("F11").Value = 10878
$ws4.Range("F20").Value = 1844
$ws4.Range("F22").Value = 589
$ws4.Range("F28").Value = 48
$ws4.Range("F29").Value = 1228
$ws4.Range("F35").Value = 1422
$ws4.Range("F38").Value = 349
$ws4.Range("F39").Value = 524
$ws4.Range("F41").Value = 358
$ws4.Range("F42").Value = 105
